# Update the test data: replace "amuthan"/"sakthivel" with "Manjeet"/"Singh"
# in row 8 (F8/G8) of the TESTDATA sheet.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("TESTDATA")

$ws2.Range("F8").Value = "Manjeet"
$ws2.Range("G8").Value = "Singh"

# Update the view/selection state: RUNMANAGER loses the selected-tab flag and
# its cursor moves to A32, while TESTDATA becomes the active/selected sheet
# with its cursor at G8.
$ws1.Range("A32").Select()
$ws2.Activate()
$ws2.Range("G8").Select()
